# Rewrite the guest-list data (rows 2-6) to reflect Devora's requested changes:
#  - move the "60107 / regular invitation" entries (testing person 1 & 2) above the family row
#  - update the family (42652) row: RSVP -> Yes, # coming -> 4, diet info -> Allergies
#  - keep the "20349 / with guest" entries (army friend, Guest) as the last two rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# headers in row 1 (unchanged):
# A invitation #, B invitation name, C guest english name, D guest hebrew name,
# E RSVP, F # invited, G # coming, H date opened, I diet info, J side, K group,
# L email, M url

$data = @(
    @("60107","regular invitation","testing person 1","טסטינג פרסון 1","Maybe",1,0," ","","Groom","Work","1@1.com","avichaidevora.com/invitation/60107"),
    @("60107","regular invitation","testing person 2","טסטינג פרסון 2","Maybe",1,0," ","","Groom","Work","2@2.com","avichaidevora.com/invitation/60107"),
    @("42652","family","the Moskovitzes","משפחת מוסקוביץ","Yes",5,4," ","Allergies","Bride","Family","reyley1014@gmail.com","avichaidevora.com/invitation/42652"),
    @("20349","with guest","army friend","חבר מהצבא","Maybe",1,0," ","","Groom","Army","reyley1014@gmail.com","avichaidevora.com/invitation/20349"),
    @("20349","with guest","Guest","אורח/ת","Maybe",1,0," ","","Groom","Army","","avichaidevora.com/invitation/20349")
)

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M")

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = $r + 2
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $columns.Length; $c++) {
        $cellRef = $columns[$c] + $rowNum
        $ws.Range($cellRef).Value = $rowValues[$c]
    }
}
